# Weekly update: insert a new data row (row 216) into the "Choclo" sheet,
# pushing all subsequent rows down by one, and fill in the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 216; everything from old row 216 onward
# shifts down to 217..316 (dimension grows from A1:R315 to A1:R316).
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new week's data.
$ws.Cells.Item(216, 1).Value2  = 7
$ws.Cells.Item(216, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(216, 3).Value2  = "Ñuble"
$ws.Cells.Item(216, 4).Value2  = 45029
$ws.Cells.Item(216, 5).Value2  = 16
$ws.Cells.Item(216, 6).Value2  = 100112024
$ws.Cells.Item(216, 7).Value2  = "Choclo"
$ws.Cells.Item(216, 8).Value2  = "Choclero"
$ws.Cells.Item(216, 9).Value2  = "Primera"
$ws.Cells.Item(216, 10).Value2 = 10000
$ws.Cells.Item(216, 11).Value2 = 400
$ws.Cells.Item(216, 12).Value2 = 400
$ws.Cells.Item(216, 13).Value2 = 400
$ws.Cells.Item(216, 14).Value2 = "`$/unidad"
$ws.Cells.Item(216, 15).Value2 = "Región del Maule"
$ws.Cells.Item(216, 16).Value2 = 400
$ws.Cells.Item(216, 17).Value2 = 1
$ws.Cells.Item(216, 18).Value2 = "Hortaliza"
